$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status names and descriptions
$ws.Range("B2").Value = "Online"
$ws.Range("B3").Value = "Offline"
$ws.Range("C2").Value = "Currently logged in."
$ws.Range("C3").Value = "Currently not logged in."

# Update the active cell selection
$ws.Range("C6").Select()
